$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.610.75"
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = "'241.98"
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D6").Value = "'0.678"
$ws.Range("E6").Value = '  +2.47%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'52.82"
$ws.Range("E8").Value = '  -6.13%  '
$ws.Range("D9").Value = "'58.87"
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("E10").Value = '  -5.56%  '
$ws.Range("E11").Value = '  -3.32%  '
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = "'0.891"
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("D14").Value = "'14.47"
$ws.Range("E14").Value = '  -9.00%  '
$ws.Range("D15").Value = "'2.369.64"
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("E16").Value = '  -4.68%  '
$ws.Range("D17").Value = "'2.058.85"
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("D18").Value = "'36.530.28"
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").Value = "'16.41"
$ws.Range("E19").Value = '  -12.60%  '
$ws.Range("D20").Value = "'71.73"
$ws.Range("E20").Value = '  -3.76%  '
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("E22").Value = '  -1.92%  '
$ws.Range("D23").Value = "'236.08"
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = '  -4.63%  '
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").Value = "'9.25"
$ws.Range("E27").Value = '  -2.87%  '
$ws.Range("D28").Value = "'163.26"
$ws.Range("E28").Value = '  -4.65%  '
$ws.Range("D29").Value = "'20.41"
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").Value = "'5.08"
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("E32").Value = '  -3.21%  '
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("E34").Value = '  -3.69%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = "'2.29"
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("D37").Value = "'1.83"
$ws.Range("E37").Value = '  -2.20%  '
$ws.Range("D38").Value = "'0.0817"
$ws.Range("E38").Value = '  -6.25%  '
$ws.Range("E39").Value = '  -5.72%  '
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").Value = "'4.87"
$ws.Range("E40").Value = '  -5.16%  '
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").Value = "'2.88"
$ws.Range("E41").Value = '  -6.15%  '
$ws.Range("D42").Value = "'0.0215"
$ws.Range("E42").Value = '  -3.08%  '
$ws.Range("E43").Value = '  -2.12%  '
$ws.Range("D44").Value = "'0.0932"
$ws.Range("E44").Value = '  -6.29%  '
$ws.Range("D45").Value = "'94.04"
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("D46").Value = "'1.387.50"
$ws.Range("E46").Value = '  +8.22%  '
$ws.Range("D47").Value = "'15.58"
$ws.Range("E47").Value = '  -9.19%  '
$ws.Range("D48").Value = "'7.31"
$ws.Range("E48").Value = '  +7.82%  '
$ws.Range("D49").Value = "'2.31"
$ws.Range("E49").Value = '  -2.62%  '
$ws.Range("D50").Value = "'2.85"
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = "'2.254.87"
$ws.Range("E51").Value = '  +1.32%  '
